$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 16362.714
$ws.Range("I9").Value = 20108
$ws.Range("K9").Value = 20108
$ws.Range("M9").Value = -19939

# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 262.54544
$ws.Range("I12").Value = 305.57144
$ws.Range("K12").Value = 305.57144
$ws.Range("M12").Value = -135.57144

# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 2168.25
$ws.Range("I19").Value = 2224.3333
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 2224.3333
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -2049.3333
$ws.Range("N19").Value = -2350

# Row 53 (Leve Item ID 5479)
$ws.Range("H53").Value = 694.3889
$ws.Range("I53").Value = 504.36365
$ws.Range("K53").Value = 504.36365
$ws.Range("M53").Value = 132.63635

# Row 88 (Leve Item ID 12608)
$ws.Range("H88").Value = 3301.1667
$ws.Range("I88").Value = 2903
$ws.Range("J88").Value = 3380.8
$ws.Range("K88").Value = 2903
$ws.Range("L88").Value = 3380.8
$ws.Range("M88").Value = -2497
$ws.Range("N88").Value = -4192.8

# Row 91 (Leve Item ID 12608)
$ws.Range("H91").Value = 3301.1667
$ws.Range("I91").Value = 2903
$ws.Range("J91").Value = 3380.8
$ws.Range("K91").Value = 2903
$ws.Range("L91").Value = 3380.8
$ws.Range("M91").Value = -1499
$ws.Range("N91").Value = -6188.8

# Row 131 (Leve Item ID 36108)
$ws.Range("H131").Value = 3647.5
$ws.Range("J131").Value = 3900
$ws.Range("L131").Value = 11700
$ws.Range("N131").Value = -21780

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 14884.566
$ws.Range("J137").Value = 7637.231
$ws.Range("L137").Value = 22911.693
$ws.Range("N137").Value = -28011.693

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 4002.4075
$ws.Range("I2").Value = 2422.8333
$ws.Range("K2").Value = 2422.8333
$ws.Range("M2").Value = -2309.8333

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3491.8447
$ws.Range("I32").Value = 3202.228
$ws.Range("K32").Value = 3202.228
$ws.Range("M32").Value = -2915.228

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 3811.2307
$ws.Range("I45").Value = 3433
$ws.Range("J45").Value = 4662.25
$ws.Range("K45").Value = 3433
$ws.Range("L45").Value = 4662.25
$ws.Range("M45").Value = -3056
$ws.Range("N45").Value = -5416.25

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 4002.4075
$ws.Range("I116").Value = 2422.8333
$ws.Range("K116").Value = 2422.8333
$ws.Range("M116").Value = -128.8332999999998

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2517.7368
$ws.Range("I122").Value = 2026.0769
$ws.Range("J122").Value = 3583
$ws.Range("K122").Value = 6078.2307
$ws.Range("L122").Value = 10749
$ws.Range("M122").Value = -3628.2307
$ws.Range("N122").Value = -15649

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 4002.4075
$ws.Range("I3").Value = 2422.8333
$ws.Range("K3").Value = 2422.8333
$ws.Range("M3").Value = -2308.8333

# Row 50 (Leve Item ID 27159)
$ws.Range("H50").Value = 63220.668
$ws.Range("J50").Value = 63220.668
$ws.Range("L50").Value = 63220.668
$ws.Range("N50").Value = -64368.668

# Row 75 (Leve Item ID 11872)
$ws.Range("H75").Value = 44664.5
$ws.Range("I75").Value = 7662.6665
$ws.Range("J75").Value = 81666.336
$ws.Range("K75").Value = 7662.6665
$ws.Range("L75").Value = 81666.336
$ws.Range("M75").Value = -6726.6665
$ws.Range("N75").Value = -83538.336

# Row 78 (Leve Item ID 11872)
$ws.Range("H78").Value = 44664.5
$ws.Range("I78").Value = 7662.6665
$ws.Range("J78").Value = 81666.336
$ws.Range("K78").Value = 22987.9995
$ws.Range("L78").Value = 244999.008
$ws.Range("M78").Value = -18307.9995
$ws.Range("N78").Value = -254359.008

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 112700
$ws.Range("I86").Value = 50400
$ws.Range("J86").Value = 175000
$ws.Range("K86").Value = 50400
$ws.Range("L86").Value = 175000
$ws.Range("M86").Value = -49277
$ws.Range("N86").Value = -177246

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 112700
$ws.Range("I89").Value = 50400
$ws.Range("J89").Value = 175000
$ws.Range("K89").Value = 252000
$ws.Range("L89").Value = 875000
$ws.Range("M89").Value = -246384
$ws.Range("N89").Value = -886232

# Row 109 (Leve Item ID 27096)
$ws.Range("H109").Value = 89791.39999999999
$ws.Range("J109").Value = 89791.39999999999
$ws.Range("L109").Value = 89791.39999999999
$ws.Range("N109").Value = -92565.39999999999

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3837.0557
$ws.Range("I31").Value = 2032.5454
$ws.Range("J31").Value = 6672.7144
$ws.Range("K31").Value = 2032.5454
$ws.Range("L31").Value = 6672.7144
$ws.Range("M31").Value = -1737.5454
$ws.Range("N31").Value = -7262.7144

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3837.0557
$ws.Range("I34").Value = 2032.5454
$ws.Range("J34").Value = 6672.7144
$ws.Range("K34").Value = 2032.5454
$ws.Range("L34").Value = 6672.7144
$ws.Range("M34").Value = -1830.5454
$ws.Range("N34").Value = -7076.7144

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 15280428
$ws.Range("I58").Value = 2564.2144
$ws.Range("K58").Value = 2564.2144
$ws.Range("M58").Value = -2361.2144

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 4059.111
$ws.Range("J134").Value = 5999.5
$ws.Range("L134").Value = 17998.5
$ws.Range("N134").Value = -23068.5

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 15280428
$ws.Range("I136").Value = 2564.2144
$ws.Range("K136").Value = 7692.6432
$ws.Range("M136").Value = -5142.6432

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 2140.182
$ws.Range("I122").Value = 801
$ws.Range("J122").Value = 2437.7778
$ws.Range("K122").Value = 7209
$ws.Range("L122").Value = 21940.0002
$ws.Range("M122").Value = -4759
$ws.Range("N122").Value = -26840.0002

# Row 134 (Leve Item ID 44074)
$ws.Range("H134").Value = 1268.8077
$ws.Range("I134").Value = 958
$ws.Range("J134").Value = 4998.5
$ws.Range("K134").Value = 2874
$ws.Range("L134").Value = 14995.5
$ws.Range("M134").Value = 2196
$ws.Range("N134").Value = -25135.5

# Row 136 (Leve Item ID 44093)
$ws.Range("H136").Value = 1554.2858
$ws.Range("I136").Value = 1554.2858
$ws.Range("K136").Value = 4662.857400000001
$ws.Range("M136").Value = 437.1425999999992

# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 2508.8333
$ws.Range("I137").Value = 2404
$ws.Range("K137").Value = 7212
$ws.Range("M137").Value = -2112

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 45 (Leve Item ID 27225)
$ws.Range("H45").Value = 37799.8
$ws.Range("I45").Value = 34999.668
$ws.Range("J45").Value = 42000
$ws.Range("K45").Value = 34999.668
$ws.Range("L45").Value = 42000
$ws.Range("M45").Value = -34440.668
$ws.Range("N45").Value = -43118

# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 9021.888999999999
$ws.Range("I70").Value = 4998.75
$ws.Range("J70").Value = 12240.4
$ws.Range("K70").Value = 4998.75
$ws.Range("L70").Value = 12240.4
$ws.Range("M70").Value = -4728.75
$ws.Range("N70").Value = -12780.4

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 9021.888999999999
$ws.Range("I73").Value = 4998.75
$ws.Range("J73").Value = 12240.4
$ws.Range("K73").Value = 4998.75
$ws.Range("L73").Value = 12240.4
$ws.Range("M73").Value = -4062.75
$ws.Range("N73").Value = -14112.4

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 2893.1333
$ws.Range("I122").Value = 1977.6666
$ws.Range("J122").Value = 4266.3335
$ws.Range("K122").Value = 5932.9998
$ws.Range("L122").Value = 12799.0005
$ws.Range("M122").Value = -3482.9998
$ws.Range("N122").Value = -17699.0005

# Row 140 (Leve Item ID 42458)
$ws.Range("H140").Value = 119271.875
$ws.Range("J140").Value = 119271.875
$ws.Range("L140").Value = 119271.875
$ws.Range("N140").Value = -129631.875

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 950
$ws.Range("J16").Value = 531
$ws.Range("L16").Value = 531
$ws.Range("N16").Value = -871

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2208.5881
$ws.Range("I46").Value = 954.8
$ws.Range("K46").Value = 954.8
$ws.Range("M46").Value = -766.8

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 4181.1
$ws.Range("I122").Value = 3962.3044
$ws.Range("K122").Value = 11886.9132
$ws.Range("M122").Value = -9436.913199999999

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2829.2222
$ws.Range("I132").Value = 2267.4443
$ws.Range("K132").Value = 6802.3329
$ws.Range("M132").Value = -4272.3329

# Row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 59680
$ws.Range("J133").Value = 59680
$ws.Range("L133").Value = 59680
$ws.Range("N133").Value = -64740

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 42 (Leve Item ID 3372)
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 3680.0952
$ws.Range("I122").Value = 2654.611
$ws.Range("J122").Value = 9833
$ws.Range("K122").Value = 7963.833
$ws.Range("L122").Value = 29499
$ws.Range("M122").Value = -5513.833
$ws.Range("N122").Value = -34399

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 3366.6365
$ws.Range("I132").Value = 3124.138
$ws.Range("K132").Value = 9372.414000000001
$ws.Range("M132").Value = -6842.414000000001

# Row 133 (Leve Item ID 41869)
$ws.Range("H133").Value = 80385
$ws.Range("J133").Value = 80385
$ws.Range("L133").Value = 80385
$ws.Range("N133").Value = -90505
